$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "tous les liens dans le footer" -> "robots,txt" entry (now with its explanation + x) ---
$ws.Range("B6").Value = "robots,txt"
$ws.Range("C6").Value = "Le fichier robots.txt ainsi que la balise meta robots est utilisé par les robots de recherches pour savoir quelle page indexer"
$ws.Range("G6").Value = "x"

# --- Row 7: becomes "meta canonical", explanation + x removed ---
$ws.Range("B7").Value = "meta canonical"
$ws.Range("C7").Clear()
$ws.Range("G7").Clear()

# --- Row 8: becomes Accessibilité / balise lang non présente ---
$ws.Range("A8").Value = "Accessibilité"
$ws.Range("B8").Value = "balise lang non présente"
$ws.Range("C8").Value = "Cette balise est utile pour les lecteurs d’écrans"
$ws.Range("D8").Value = "ajouter lang=’fr’ dans la balise html"
$ws.Range("F8").Value = "lighthouse"
$ws.Range("G8").Value = "x"

# --- Row 9: becomes (SEO ou accessiblité ?) / Cache pour les photos ---
$ws.Range("A9").Value = "(SEO ou accessiblité ?)"
$ws.Range("B9").Value = "Cache pour les photos"
$ws.Range("C9").Clear()
$ws.Range("D9").Clear()
$ws.Range("F9").Clear()
$ws.Range("G9").Clear()

# --- Row 10: becomes Vulnérabilité / Jquery upgrade note, no formula anymore ---
$ws.Range("A10").Value = "Vulnérabilité"
$ws.Range("B10").Value = "Jquery 2,1,0 à jquery 3,5,1"
$ws.Range("E10").Clear()

# --- Row 11: becomes Accessibilité / label sur social ---
$ws.Range("A11").Value = "Accessibilité"
$ws.Range("B11").Value = "label sur social"
$ws.Range("C11").Value = "Les liens de réseaux sociaux n’ont pas de label, ce qui bloque le lecteur d’écran"
$ws.Range("D11").Value = "ajouter aria-label=’’ dans les balises liens"
$ws.Range("F11").Value = "lighthouse"
$ws.Range("G11").Value = "x"

# --- Row 12: becomes (SEO ou accessiblité ?) / minifier le css ---
$ws.Range("A12").Value = "(SEO ou accessiblité ?)"
$ws.Range("B12").Value = "minifier le css"
$ws.Range("C12").Value = "Ces fichiers prennent de la place, les minifiers permet de gagner en rapidité"
$ws.Range("D12").Clear()
$ws.Range("F12").Clear()
$ws.Range("G12").Clear()

# --- Row 13: becomes Bug dans bloc.js ---
$ws.Range("B13").Value = "Bug dans bloc.js"
$ws.Range("C13").Clear()

# --- Row 14: becomes SEO / Certains texte font 1px ---
$ws.Range("A14").Value = "SEO"
$ws.Range("B14").Value = "Certains texte font 1px"
$ws.Range("C14").Value = "Le fait d’avoir du texte de très petite taille ou de la meme couleur que le fond peu etre considéré comme du hack"

# --- Row 15: new footer-links entry with a normal (non-bold-ish SEO-category) style and a reference link ---
$ws.Range("A15").Font.Name = "Arial"
$ws.Range("A15").Font.Color = 0
$ws.Range("B15").Value = "Beaucoups de liens dans le footer"
$ws.Range("C15").Clear()
$ws.Range("F15").Value = "https://www.rocktherankings.com/footer-links-seo/"

# --- Update selection to match the new cursor position ---
$ws.Range("C17").Select() | Out-Null
